$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Row 2 : the former "nothing" placeholder text becomes the real
#         expected-result text for the first (already existing) case.
# -----------------------------------------------------------------
$ws.Range("E2").Value = "BMI = 13.9 kg/m2   (Severe thinness)"

# -----------------------------------------------------------------
# Row 3 : new height/weight, rich-text expected result with a
#         superscript "2" and a red/bold "Severe thinness" phrase.
# -----------------------------------------------------------------
$ws.Range("C3").Value = 180
$ws.Range("D3").Value = 51.516

$e3 = $ws.Range("E3")
$e3.Value = "BMI = 15.9 kg/m2   (Severe thinness)"

# Whole-cell base font: bold, 14pt, black, Arial. Building it through a
# transient named style (instead of mutating $e3.Font property-by-property)
# avoids leaving a trail of half-applied intermediate fonts in the
# styles part.
$tmpStyle = $wb.Styles.Add("TmpBmiStyle")
$tmpStyle.Font.Bold = $true
$tmpStyle.Font.Size = 14
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Color = 0
$e3.Style = "TmpBmiStyle"
$wb.Styles.Item("TmpBmiStyle").Delete()

# "2" (the exponent in kg/m2) -> superscript
$e3.Characters(16, 1).Font.Superscript = $true

# "   (" -> not bold
$e3.Characters(17, 4).Font.Bold = $false

# "Severe thinness" -> bold + red
$sev = $e3.Characters(21, 15)
$sev.Font.Color = 394937
$sev.Font.Bold = $true

# ")" -> not bold
$e3.Characters(36, 1).Font.Bold = $false

# Row 3 is taller to fit the bigger font, and new column E is wider.
$ws.Rows.Item(3).RowHeight = 21
$ws.Columns.Item(5).ColumnWidth = 35.6

# -----------------------------------------------------------------
# Row 4 : new height/weight/expected-result.
# -----------------------------------------------------------------
$ws.Range("C4").Value = 175
$ws.Range("D4").Value = 49
$ws.Range("E4").Value = "BMI = 16 kg/m2   (Moderate thinness)"

# -----------------------------------------------------------------
# Rows 5-16 : brand-new BMI test cases.
# -----------------------------------------------------------------
$newRows = @(
    @{ Row = 5;  Age = 40; Height = 170; Weight = 49.13;  Result = "BMI = 17 kg/m2   (Moderate thinness)" },
    @{ Row = 6;  Age = 40; Height = 190; Weight = 61.731; Result = "BMI = 17.1 kg/m2   (Mild thinness)" },
    @{ Row = 7;  Age = 40; Height = 190; Weight = 66.785; Result = "BMI = 18.5 kg/m2   (Mild thinness)" },
    @{ Row = 8;  Age = 40; Height = 190; Weight = 67.146; Result = "BMI = 18.6 kg/m2   (Normal)" },
    @{ Row = 9;  Age = 40; Height = 190; Weight = 90.25;  Result = "BMI = 25 kg/m2   (Normal)" },
    @{ Row = 10; Age = 40; Height = 190; Weight = 90.611; Result = "BMI = 25.1 kg/m2   (Overweight)" },
    @{ Row = 11; Age = 40; Height = 200; Weight = 120;    Result = "BMI = 30 kg/m2   (Overweight)" },
    @{ Row = 12; Age = 40; Height = 200; Weight = 120.4;  Result = "BMI = 30.1 kg/m2   (Obese Class I)" },
    @{ Row = 13; Age = 40; Height = 200; Weight = 140;    Result = "BMI = 35 kg/m2   (Obese Class I)" },
    @{ Row = 14; Age = 40; Height = 200; Weight = 140.4;  Result = "BMI = 35.1 kg/m2   (Mild thinness)" },
    @{ Row = 15; Age = 40; Height = 200; Weight = 160;    Result = "BMI = 40 kg/m2   (Mild thinness)" },
    @{ Row = 16; Age = 40; Height = 200; Weight = 160.4;  Result = "BMI = 40.1 kg/m2   (Mild thinness)" }
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.Age
    $ws.Range("B" + $r.Row).Value = "Male"
    $ws.Range("C" + $r.Row).Value = $r.Height
    $ws.Range("D" + $r.Row).Value = $r.Weight
    $ws.Range("E" + $r.Row).Value = $r.Result
}

# -----------------------------------------------------------------
# Sheet-level tweaks: selection, print orientation.
# -----------------------------------------------------------------
$ws.Range("D17").Select() | Out-Null
$ws.PageSetup.Orientation = 1
